$d = $word.ActiveDocument

# 1. Replace the subtitle text "RequestSolved!" with the new project title.
$d.Content.Find.Execute("RequestSolved!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "E-commerce de joias e itens de artesanato", 2)

# 2. Remove the "Tema" summary table that followed the subtitle.
$d.Tables.Item(1).Delete()

# 3. Remove the two now-orphaned empty paragraphs that used to sit between the
#    table and the "N01:" requirements list, without disturbing the subtitle
#    paragraph's own ending mark (which carries the "Subtitulo" style) or the
#    "N01:" paragraph that follows.
$rAfterSubtitle = $d.Content
$rAfterSubtitle.Find.Execute("E-commerce de joias e itens de artesanato", $true, $false, $false,
                              $false, $false, $true, 1, $false, "", 0)
$afterSubtitle = $rAfterSubtitle.End

$rBeforeN01 = $d.Content
$rBeforeN01.Find.Execute("N01:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$beforeN01 = $rBeforeN01.Start

# Skip past the subtitle paragraph's own paragraph mark (afterSubtitle) and
# delete the remaining paragraph marks one at a time; each successful delete
# shrinks the gap by one character.
$pos = $afterSubtitle + 1
while ($pos -lt $beforeN01) {
    $charRange = $d.Range($pos, $pos + 1)
    $charRange.Delete()
    $beforeN01 = $beforeN01 - 1
}
